# Remove the debug yellow (and any other) highlighting that was left in the
# final document templates. Track Changes must be off, otherwise Word would
# just wrap the change in <w:rPrChange> markup instead of really removing it.
$d = $word.ActiveDocument
$d.TrackRevisions = $false

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Highlight = $true
$find.Text = ""
$find.Replacement.Highlight = $false
$find.Replacement.Text = ""
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $true, $null, 2) | Out-Null
